# "đưa lên mạch I2C" — the "Rơ le" module (row 16) now shares the I2C board,
# so its "thiếu Footprint relay G6S-2-12VDC" note cell (G16) is marked with
# the same green, centered "resolved" styling used elsewhere in column G.
#
# Also: the stray leftover note "mung lung như một trò đùa :(((" is removed
# from the four rows that still referenced it (F28, F29, F30, F36), and the
# active selection left on the sheet moves from F16 to F36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# G16: give it the green "done" fill + centered alignment (same look as
# G4/G5/G23/G24/etc.), matching the border it already has.
$g16 = $ws.Range("G16")
$g16.HorizontalAlignment = $xlCenter
$g16.VerticalAlignment = $xlCenter
$g16.Interior.Pattern = 1
$g16.Interior.Color = 5296274

# Clear the obsolete "mung lung như một trò đùa :(((" note wherever it still
# lingers, leaving the cells blank (formatting untouched).
$ws.Range("F28").ClearContents()
$ws.Range("F29").ClearContents()
$ws.Range("F30").ClearContents()
$ws.Range("F36").ClearContents()

# Leave the sheet's selection where the author last left it.
[void]$ws.Range("F36").Select()
